# Added new ss capture method and changed Driver to static
$wb = $excel.ActiveWorkbook

# Work on the "BDD Scenario" sheet: add the new Feature/Scenario rows
$ws = $wb.Worksheets.Item("BDD Scenario")

$ws.Range("B12").Value = "Feature: Place Order"
$ws.Range("B14").Value = "Scenario: Place order on SauceDemo"

# Make "BDD Scenario" the active/selected sheet and set the new selection
$ws.Activate()
$ws.Range("B17").Select()
